# Generate Report for Handoff
# The localization CI run regenerated its report against a newly-minted
# source GUID, so every cell that echoed the old GUID-named file needs to
# reflect the new one, and the handoff timestamps need to catch up to the
# moment this run actually finished.

$wb = $excel.ActiveWorkbook

$oldGuid = "df47ccd6-2f99-4912-bd61-a970e014bd1b"
$newGuid = "94f05208-5e43-45fc-adab-6534918ae77f"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$newPathAndName = "e2e\$newGuid.md"
$wsOverview.Range("B2").Value = $newPathAndName
# Re-point the hyperlink so its display text tracks the new file name too
# (the underlying target address is left as-is, same as upstream).
$overviewLinkAddr = $wsOverview.Range("B2").Hyperlinks.Item(1).Address
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewLinkAddr, "", "", $newPathAndName)

$wsOverview.Range("G2").Value = "2016-08-29 10:58:46"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$newFileName = "$newGuid.md"
$wsZhCn.Range("A2").Value = $newFileName
$zhCnLinkAddr = $wsZhCn.Range("A2").Hyperlinks.Item(1).Address
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhCnLinkAddr, "", "", $newFileName)

$wsZhCn.Range("G2").Value = "$newGuid.25cd6fb6bdd7cd55cb3a46d7b0a6421a09a613e8.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-29 10:58:42"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newFileName
$deDeLinkAddr = $wsDeDe.Range("A2").Hyperlinks.Item(1).Address
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $deDeLinkAddr, "", "", $newFileName)

$wsDeDe.Range("G2").Value = "$newGuid.25cd6fb6bdd7cd55cb3a46d7b0a6421a09a613e8.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-29 10:58:46"
